# Distributed System Report.docx - edit script
#
# 1) Expand the "Client communicates with Master and Slaves" sentence.
# 2) Add a numbered ("Writing operation" / "Reading operation") list and
#    the two explanatory paragraphs describing the write & read protocols.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Rewrite the closing sentence of the intro paragraph.
# ---------------------------------------------------------------------
$oldSentence = "The Client performs two operations: write and read. During these operations, Client communicates with Master and Slaves"
$newSentence = "The Client performs two operations: write and read. During these operations, Client communicates with the Master for the Metadata and the Slaves for transferring the data."

$found = $d.Content.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# ---------------------------------------------------------------------
# 2. Locate the paragraph that now holds the rewritten sentence (still
#    paragraph 3 - "Distributed System Report" / "Client" / this one)
#    and append the new content right after it, before the bookmark's
#    paragraph mark.
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs(3)
$introEnd = $introPara.Range
$introEnd.Collapse(0)

# --- "Writing operation" (numbered list item 1) ---
$introEnd.InsertParagraphAfter()
$writingHeading = $d.Paragraphs(4)
$writingHeading.Range.Text = "Writing operation"
$writingHeading.Style = "List Paragraph"
$writingHeading.Range.ListFormat.ApplyNumberDefault()

# --- Paragraph describing the write operation ---
$writingHeadingEnd = $writingHeading.Range
$writingHeadingEnd.Collapse(0)
$writingHeadingEnd.InsertParagraphAfter()
$writeBody = $d.Paragraphs(5)
$writeBody.Range.Text = "A write operation is executed in 3 steps. First, the Client makes a request to the Master to notify about this operation. The original file" + [char]8217 + "s size is also included in this request for the decoding step in the future. The Master sends back the list of live Slaves where data will be written to. In the case of Hierarchical code, the number of Slaves in the list is 7. Each of these Slaves is associated to an encoded part which will be stored on that Slave. Then, the Client encodes the file into several parts using erasure codes. After that, the Client contacts to each of Slaves on the list and transfers the corresponding encoded data. However, there may be failures in the communication with the Slaves. These failures may be caused by a communication problem or the fact that the Slave already died and the Master has not known about this (because of the maintaining process). We assume that there is only second type of failure in our experiments. In that case, the Client should ignore these errors because sooner or later the Master will know about this dead Slave. The Client does not have to acknowledge the Master about the failed Slave, which keeps the system still simple and efficient in failure handling."
$writeBody.Range.ParagraphFormat.LeftIndent = 18

# --- "Reading operation" (numbered list item 2) ---
$writeBodyEnd = $writeBody.Range
$writeBodyEnd.Collapse(0)
$writeBodyEnd.InsertParagraphAfter()
$readingHeading = $d.Paragraphs(6)
$readingHeading.Range.Text = "Reading operation"
$readingHeading.Style = "List Paragraph"
$readingHeading.Range.ListFormat.ApplyNumberDefault()

# --- Paragraph describing the read operation ---
$readingHeadingEnd = $readingHeading.Range
$readingHeadingEnd.Collapse(0)
$readingHeadingEnd.InsertParagraphAfter()
$readBody = $d.Paragraphs(7)
$readBody.Range.Text = "Similarly, a read operation includes 3 phases. At the beginning, the Client requests to the Master for the Metadata of the file. The Master returns all the live Slaves containing encoded parts as well as the original file" + [char]8217 + "s size. Based on this list, the Client determines which Slaves it should contact to retrieve data and reconstruct the file. The Client can communicate to each Slave sequentially or concurrently for a better performance. Similar to the writing operation, there may be a chance of failures because of 2 above reasons. In both cases, the Client has to recalculate to decide one or more other parts needed to recover data; therefore, more communications are required. With these collected pieces of information, the decoding process is executed by the Client. If the original file cannot be reconstructed in some fault scenarios, the read operation is unsuccessful. "
$readBody.Range.ParagraphFormat.LeftIndent = 18

# ---------------------------------------------------------------------
# 3. Tidy up the "List Paragraph" style definition that Word mints on
#    first use so it matches the built-in gallery definition.
# ---------------------------------------------------------------------
$listStyle = $d.Styles("List Paragraph")
$listStyle.Priority = 34
$listStyle.ParagraphFormat.LeftIndent = 36

Write-Output "done"
